# Update the "想去人数" (Column F) values on the "展览" and "全部类型"
# worksheets to reflect the latest scrape, as described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F (same update applies to
# both the "展览" sheet and the "全部类型" sheet, which mirror each other).
$updates = @{
    3  = 3055
    4  = 221
    5  = 117
    6  = 195
    7  = 1651
    8  = 1617
    15 = 225
    16 = 233
    20 = 41
    22 = 361
    23 = 168
    26 = 2044
    29 = 17
    30 = 185
    34 = 5
    35 = 495
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
